$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($sheet, $addr, $text)
    $rng = $sheet.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-CellText $ws "D2" "308.31"
Set-CellText $ws "E2" "0.24%"
Set-CellText $ws "D3" "40.85"
Set-CellText $ws "E3" "2.41%"
Set-CellText $ws "D4" "5.108"
Set-CellText $ws "E4" "-0.21%"
Set-CellText $ws "D5" "0.07624"
Set-CellText $ws "E5" "-1.31%"
Set-CellText $ws "D6" "4.254"
Set-CellText $ws "E6" "0.58%"
Set-CellText $ws "D7" "1.606"
Set-CellText $ws "E7" "-0.04%"
Set-CellText $ws "E8" "2.09%"
Set-CellText $ws "E9" "0.93%"
Set-CellText $ws "D10" "0.1122"
Set-CellText $ws "E10" "11.71%"
Set-CellText $ws "D11" "0.1794"
Set-CellText $ws "E11" "3.50%"
Set-CellText $ws "D12" "0.09178"
Set-CellText $ws "E12" "1.71%"
Set-CellText $ws "D13" "0.04174"
Set-CellText $ws "E13" "-6.23%"
Set-CellText $ws "D14" "0.1052"
Set-CellText $ws "E14" "-0.11%"
Set-CellText $ws "D15" "0.001256"
Set-CellText $ws "E15" "-1.02%"
Set-CellText $ws "D16" "0.005693"
Set-CellText $ws "E16" "-1.94%"
Set-CellText $ws "E17" "-0.13%"
Set-CellText $ws "E18" "-0.73%"
Set-CellText $ws "D19" "6.656"
Set-CellText $ws "E19" "-5.59%"
Set-CellText $ws "D20" "0.1360"
Set-CellText $ws "E20" "0.93%"
Set-CellText $ws "D21" "0.2802"
Set-CellText $ws "E22" "-1.58%"
Set-CellText $ws "D23" "0.001246"
Set-CellText $ws "E23" "3.13%"
Set-CellText $ws "D24" "0.004103"
Set-CellText $ws "E24" "1.05%"
Set-CellText $ws "E25" "-0.09%"
Set-CellText $ws "D38" "0.02394"
Set-CellText $ws "E38" "2.26%"
Set-CellText $ws "D39" "0.05186"
Set-CellText $ws "E39" "-0.15%"
Set-CellText $ws "D40" "0.007787"
Set-CellText $ws "E40" "-1.64%"
Set-CellText $ws "D41" "0.1298"
Set-CellText $ws "E42" "12.99%"
Set-CellText $ws "E43" "-0.07%"
Set-CellText $ws "D44" "0.007718"
Set-CellText $ws "E44" "-6.17%"
Set-CellText $ws "D45" "0.3079"
Set-CellText $ws "E45" "-7.52%"
Set-CellText $ws "D46" "0.00006986"
Set-CellText $ws "E46" "7.26%"
Set-CellText $ws "E47" "-0.11%"
Set-CellText $ws "D48" "0.04673"
Set-CellText $ws "E48" "1,256.41%"
Set-CellText $ws "E50" "-0.11%"
Set-CellText $ws "E51" "-0.11%"
